$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 4 de Julio de 2020 a las 18:55"

# Row 4
$ws.Cells.Item(4, 2).Value = 2912169
$ws.Cells.Item(4, 3).Value = 21581
$ws.Cells.Item(4, 4).Value = 1237836
$ws.Cells.Item(4, 5).Value = 1542137
$ws.Cells.Item(4, 7).Value = 95
$ws.Cells.Item(4, 8).Value = 132196

# Row 5
$ws.Cells.Item(5, 2).Value = 1550176
$ws.Cells.Item(5, 3).Value = 6835
$ws.Cells.Item(5, 5).Value = 508152
$ws.Cells.Item(5, 7).Value = 155
$ws.Cells.Item(5, 8).Value = 63409

# Row 7
$ws.Cells.Item(7, 2).Value = 672644
$ws.Cells.Item(7, 3).Value = 22755
$ws.Cells.Item(7, 4).Value = 408625
$ws.Cells.Item(7, 5).Value = 244740
$ws.Cells.Item(7, 7).Value = 610
$ws.Cells.Item(7, 8).Value = 19279

# Row 11
$ws.Cells.Item(11, 2).Value = 284900
$ws.Cells.Item(11, 3).Value = 624
$ws.Cells.Item(11, 7).Value = 67
$ws.Cells.Item(11, 8).Value = 44198

# Row 13
$ws.Cells.Item(13, 2).Value = 241419
$ws.Cells.Item(13, 3).Value = 235
$ws.Cells.Item(13, 4).Value = 191944
$ws.Cells.Item(13, 5).Value = 14621
$ws.Cells.Item(13, 7).Value = 21
$ws.Cells.Item(13, 8).Value = 34854

# Row 18
$ws.Cells.Item(18, 2).Value = 197250
$ws.Cells.Item(18, 3).Value = 250
$ws.Cells.Item(18, 5).Value = 6876

# Row 23
$ws.Cells.Item(23, 2).Value = 105211
$ws.Cells.Item(23, 3).Value = 120
$ws.Cells.Item(23, 4).Value = 68868
$ws.Cells.Item(23, 5).Value = 27675
$ws.Cells.Item(23, 7).Value = 5
$ws.Cells.Item(23, 8).Value = 8668

# Row 40
$ws.Cells.Item(40, 4).Value = 40117
$ws.Cells.Item(40, 5).Value = 4521

# Row 54
$ws.Cells.Item(54, 2).Value = 25509
$ws.Cells.Item(54, 3).Value = 11
$ws.Cells.Item(54, 5).Value = 404
$ws.Cells.Item(54, 7).Value = 1
$ws.Cells.Item(54, 8).Value = 1741

# Row 57
$ws.Cells.Item(57, 1).Value = "Azerbaiyan"
$ws.Cells.Item(57, 2).Value = 19801
$ws.Cells.Item(57, 3).Value = 534
$ws.Cells.Item(57, 4).Value = 11291
$ws.Cells.Item(57, 5).Value = 8269
$ws.Cells.Item(57, 7).Value = 6
$ws.Cells.Item(57, 8).Value = 241

# Row 58
$ws.Cells.Item(58, 1).Value = "Ghana"
$ws.Cells.Item(58, 2).Value = 19388
$ws.Cells.Item(58, 4).Value = 14330
$ws.Cells.Item(58, 5).Value = 4941
$ws.Cells.Item(58, 8).Value = 117

# Row 63
$ws.Cells.Item(63, 1).Value = "Argelia"
$ws.Cells.Item(63, 2).Value = 15500
$ws.Cells.Item(63, 3).Value = 430
$ws.Cells.Item(63, 4).Value = 11181
$ws.Cells.Item(63, 5).Value = 3373
$ws.Cells.Item(63, 7).Value = 9
$ws.Cells.Item(63, 8).Value = 946

# Row 64
$ws.Cells.Item(64, 1).Value = "Nepal"
$ws.Cells.Item(64, 2).Value = 15491
$ws.Cells.Item(64, 3).Value = 232
$ws.Cells.Item(64, 4).Value = 6415
$ws.Cells.Item(64, 5).Value = 9042
$ws.Cells.Item(64, 7).Value = 2
$ws.Cells.Item(64, 8).Value = 34

# Row 72
$ws.Cells.Item(72, 2).Value = 9611
$ws.Cells.Item(72, 3).Value = 215
$ws.Cells.Item(72, 4).Value = 6306
$ws.Cells.Item(72, 5).Value = 3274
$ws.Cells.Item(72, 7).Value = 2
$ws.Cells.Item(72, 8).Value = 31

# Row 73
$ws.Cells.Item(73, 2).Value = 8926
$ws.Cells.Item(73, 3).Value = 5
$ws.Cells.Item(73, 5).Value = 537

# Row 82
$ws.Cells.Item(82, 2).Value = 6932
$ws.Cells.Item(82, 3).Value = 145
$ws.Cells.Item(82, 4).Value = 2987
$ws.Cells.Item(82, 5).Value = 3611
$ws.Cells.Item(82, 7).Value = 6
$ws.Cells.Item(82, 8).Value = 334

# Row 123
$ws.Cells.Item(123, 1).Value = "Sierra Leona"
$ws.Cells.Item(123, 2).Value = 1533
$ws.Cells.Item(123, 3).Value = 9
$ws.Cells.Item(123, 4).Value = 1051
$ws.Cells.Item(123, 5).Value = 420
$ws.Cells.Item(123, 8).Value = 62

# Row 124
$ws.Cells.Item(124, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(124, 2).Value = 1530
$ws.Cells.Item(124, 4).Value = 1490
$ws.Cells.Item(124, 5).Value = 18
$ws.Cells.Item(124, 8).Value = 22

# Row 130
$ws.Cells.Item(130, 2).Value = 1186
$ws.Cells.Item(130, 3).Value = 5
$ws.Cells.Item(130, 4).Value = 1046
$ws.Cells.Item(130, 5).Value = 90

# Row 131
$ws.Cells.Item(131, 2).Value = 1150
$ws.Cells.Item(131, 3).Value = 3
$ws.Cells.Item(131, 4).Value = 903
$ws.Cells.Item(131, 5).Value = 237

# Row 137
$ws.Cells.Item(137, 1).Value = "Mozambique"
$ws.Cells.Item(137, 2).Value = 969
$ws.Cells.Item(137, 3).Value = 30
$ws.Cells.Item(137, 4).Value = 256
$ws.Cells.Item(137, 5).Value = 706
$ws.Cells.Item(137, 7).Value = 1
$ws.Cells.Item(137, 8).Value = 7

# Row 138
$ws.Cells.Item(138, 1).Value = "Suazilandia"
$ws.Cells.Item(138, 2).Value = 954
$ws.Cells.Item(138, 3).Value = 45
$ws.Cells.Item(138, 4).Value = 535
$ws.Cells.Item(138, 5).Value = 406
$ws.Cells.Item(138, 8).Value = 13

# Row 139
$ws.Cells.Item(139, 1).Value = "Uruguay"
$ws.Cells.Item(139, 2).Value = 952
$ws.Cells.Item(139, 4).Value = 837
$ws.Cells.Item(139, 5).Value = 87
$ws.Cells.Item(139, 8).Value = 28

# Row 140
$ws.Cells.Item(140, 1).Value = "Georgia"
$ws.Cells.Item(140, 2).Value = 948
$ws.Cells.Item(140, 3).Value = 5
$ws.Cells.Item(140, 4).Value = 825
$ws.Cells.Item(140, 5).Value = 108
$ws.Cells.Item(140, 8).Value = 15

# Row 141
$ws.Cells.Item(141, 1).Value = "Uganda"
$ws.Cells.Item(141, 2).Value = 927
$ws.Cells.Item(141, 3).Value = 16
$ws.Cells.Item(141, 4).Value = 868
$ws.Cells.Item(141, 5).Value = 59
$ws.Cells.Item(141, 8).Value = 0

# Row 142
$ws.Cells.Item(142, 1).Value = "Libia"
$ws.Cells.Item(142, 2).Value = 918
$ws.Cells.Item(142, 4).Value = 230
$ws.Cells.Item(142, 5).Value = 661
$ws.Cells.Item(142, 8).Value = 27

# Row 147
$ws.Cells.Item(147, 1).Value = "Montenegro"
$ws.Cells.Item(147, 2).Value = 720
$ws.Cells.Item(147, 3).Value = 57
$ws.Cells.Item(147, 4).Value = 315
$ws.Cells.Item(147, 5).Value = 392

# Row 148
$ws.Cells.Item(148, 1).Value = "Santo Tome y Principe"
$ws.Cells.Item(148, 2).Value = 719
$ws.Cells.Item(148, 4).Value = 267
$ws.Cells.Item(148, 5).Value = 439

# Row 149
$ws.Cells.Item(149, 1).Value = "Crucero"
$ws.Cells.Item(149, 2).Value = 712
$ws.Cells.Item(149, 4).Value = 651
$ws.Cells.Item(149, 5).Value = 48
$ws.Cells.Item(149, 8).Value = 13

# Row 150
$ws.Cells.Item(150, 1).Value = "San Marino"
$ws.Cells.Item(150, 2).Value = 698
$ws.Cells.Item(150, 4).Value = 656
$ws.Cells.Item(150, 5).Value = 0
$ws.Cells.Item(150, 8).Value = 42

# Row 151
$ws.Cells.Item(151, 1).Value = "Malta"
$ws.Cells.Item(151, 2).Value = 672
$ws.Cells.Item(151, 4).Value = 650
$ws.Cells.Item(151, 5).Value = 13
$ws.Cells.Item(151, 8).Value = 9

# Row 152
$ws.Cells.Item(152, 1).Value = "Togo"
$ws.Cells.Item(152, 2).Value = 671
$ws.Cells.Item(152, 4).Value = 424
$ws.Cells.Item(152, 5).Value = 233
$ws.Cells.Item(152, 8).Value = 14

# Row 209
$ws.Cells.Item(209, 1).Value = "Islas Malvinas"

# Row 210
$ws.Cells.Item(210, 1).Value = "Groenlandia"
